# Apply "new cast tests added" edit to the RegressionTests sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegressionTests")

# --- Step 1: insert a single row above row 163 -------------------------
# This pushes the existing row 163 (CastInExprsDynError1) down to row 164
# and the existing row 164 (CastInExprsDynError2) down to row 165, leaving
# row 163 itself blank.
$ws.Rows("163:163").Insert()
$ws.Rows("163:163").Clear()

# --- Step 2: insert 4 new rows above the (new) row 166 -----------------
# Row 166 currently holds nothing (blank gap); this shifts everything
# from row 166 onward down by 4, opening up rows 166-169 for the new
# test entries.
$ws.Rows("166:169").Insert()

# --- Step 3: populate the 4 new rows with the new cast tests ------------
$newTests = @(
    @{ Row = 166; Name = "CastInExprsDynError3" },
    @{ Row = 167; Name = "CastInExprsDynError4" },
    @{ Row = 168; Name = "CastInExprsDynError5" },
    @{ Row = 169; Name = "CastInExprsDynError6" }
)

foreach ($t in $newTests) {
    $r = $t.Row
    $ws.Cells.Item($r, 1).Value = $t.Name
    $ws.Cells.Item($r, 2).Value = "3.4. Cast operator in expressions"
    $ws.Cells.Item($r, 3).Value = "No"
    $ws.Cells.Item($r, 4).Value = "Yes"
    $ws.Cells.Item($r, 6).Value = "Yes"
}

# --- Step 4: update view state (scroll position / selection) -----------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 131
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A169").Select()
